# "continent B coast, islands, exit"
#
# Applies:
#  - Row 3 (Continent Boscalis total): mark as done (green fill) and overwrite
#    the computed "Area (actual)" (M3) with a hand-entered value, which
#    detaches it from the K3-derived shared formula. M4 keeps computing =K4.
#  - S4: note "Microcontintent"
#  - M39: updated actual area for that island entry
#  - Rows 43/44: two new entries (Continent A SZ Islands, Continent B Sz
#    Island) filled in and marked done (green fill), matching the styling
#    already used by the neighbouring finished rows (38-42).
#  - Move the active selection to S5 (matches the saved cursor position).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$green = 5296274   # RGB(146, 208, 80) == the existing "done" fill used elsewhere in this sheet

# --- Row 3: Continent Boscalis total gets typed over with a hand value ---
$ws.Range("M3").Value = 55094497.448899999
$ws.Range("J3:M3").Interior.Color = $green

# --- M39: updated actual area value ---
$ws.Range("M39").Value = 334338.82849999995

# --- Row 43: Continent A SZ Islands ---
$ws.Range("J43").Value = "Continent A SZ Islands"
$ws.Range("M43").Value = 27742.253799999999
$ws.Range("J43:O43").Interior.Color = $green

# --- Row 44: Continent B Sz Island ---
$ws.Range("J44").Value = "Continent B Sz Island"
$ws.Range("M44").Value = 30816.3354
$ws.Range("J44:O44").Interior.Color = $green

# --- S4: add note ---
$ws.Range("S4").Value = "Microcontintent"

# --- Selection moved to S5 on save ---
$ws.Range("S5").Select()
